$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 10-16 (B column labels) to reflect the new scheme names.
$ws.Range("B10").Value = "Gaussian-Quadrature"
$ws.Range("B11").Value = "Spiral-90deg-10rot-5space"
$ws.Range("B12").Value = "Spiral-90deg-15rot-5space"
$ws.Range("B13").Value = "Spiral-90deg-10rot-3space"
$ws.Range("B14").Value = "NoRotation-tilt60deg"
$ws.Range("B15").Value = "Rotation-NoTilt"
$ws.Range("B16").Value = "Rotation-60detTilt"

# Append new rows 17-19 for the remaining averaging schemes.
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "HexGrid-90degTilt5degRes"

$ws.Range("A18").Value = 16
$ws.Range("B18").Value = "HexGrid-90degTilt22p5degRes"

$ws.Range("A19").Value = 17
$ws.Range("B19").Value = "HexGrid-60degTilt5degRes"

"C","D","E","F","G","H","I","J","K","L","M" | ForEach-Object {
    $col = $_
    $ws.Range($col + "17").Value = 1
    $ws.Range($col + "18").Value = 1
    $ws.Range($col + "19").Value = 1
}

# Copy the A column's header-style formatting (bold, bordered, centered)
# from the row above down onto the three new A-column cells.
$ws.Range("A16").Copy() | Out-Null
$ws.Range("A17:A19").PasteSpecial(-4122) | Out-Null
